$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1515
$ws.Range("I5").Value = 227
$ws.Range("J5").Value = 1998
$ws.Range("K5").Value = 227
$ws.Range("L5").Value = 1998
$ws.Range("M5").Value = -112
$ws.Range("N5").Value = -2228
$ws.Range("H17").Value = 1233.1212
$ws.Range("J17").Value = 1209.1562
$ws.Range("L17").Value = 3627.4686
$ws.Range("N17").Value = -3963.4686
$ws.Range("H19").Value = 1958.3462
$ws.Range("J19").Value = 1444
$ws.Range("L19").Value = 1444
$ws.Range("N19").Value = -1794
$ws.Range("H28").Value = 881.28
$ws.Range("I28").Value = 409.5238
$ws.Range("J28").Value = 3358
$ws.Range("K28").Value = 409.5238
$ws.Range("L28").Value = 3358
$ws.Range("M28").Value = 75.47620000000001
$ws.Range("N28").Value = -4328
$ws.Range("H40").Value = 2817.75
$ws.Range("I40").Value = 1809.4
$ws.Range("K40").Value = 1809.4
$ws.Range("M40").Value = -1634.4
$ws.Range("H51").Value = 62503844
$ws.Range("I51").Value = 125003440
$ws.Range("K51").Value = 125003440
$ws.Range("M51").Value = -125002956
$ws.Range("H62").Value = 13151.3
$ws.Range("I62").Value = 16213.571
$ws.Range("J62").Value = 6006
$ws.Range("K62").Value = 16213.571
$ws.Range("L62").Value = 6006
$ws.Range("M62").Value = -15589.571
$ws.Range("N62").Value = -7254
$ws.Range("H65").Value = 13151.3
$ws.Range("I65").Value = 16213.571
$ws.Range("J65").Value = 6006
$ws.Range("K65").Value = 81067.855
$ws.Range("L65").Value = 30030
$ws.Range("M65").Value = -77947.855
$ws.Range("N65").Value = -36270
$ws.Range("H70").Value = 1959.5385
$ws.Range("I70").Value = 2015.75
$ws.Range("K70").Value = 6047.25
$ws.Range("M70").Value = -5777.25
$ws.Range("H73").Value = 1959.5385
$ws.Range("I73").Value = 2015.75
$ws.Range("K73").Value = 6047.25
$ws.Range("M73").Value = -5111.25
$ws.Range("H80").Value = 1131.2858
$ws.Range("I80").Value = 926.6667
$ws.Range("K80").Value = 2780.0001
$ws.Range("M80").Value = -1782.0001
$ws.Range("H83").Value = 1131.2858
$ws.Range("I83").Value = 926.6667
$ws.Range("K83").Value = 8340.0003
$ws.Range("M83").Value = -3348.0003
$ws.Range("H97").Value = 1420.1765
$ws.Range("J97").Value = 1458.9375
$ws.Range("L97").Value = 4376.8125
$ws.Range("N97").Value = -5368.8125
$ws.Range("H98").Value = 917.94116
$ws.Range("I98").Value = 760.6
$ws.Range("K98").Value = 760.6
$ws.Range("M98").Value = 737.4
$ws.Range("H101").Value = 437
$ws.Range("I101").Value = 259.83334
$ws.Range("J101").Value = 1500
$ws.Range("K101").Value = 779.5000200000001
$ws.Range("L101").Value = 4500
$ws.Range("M101").Value = 842.4999799999999
$ws.Range("N101").Value = -7744
$ws.Range("H105").Value = 73000
$ws.Range("J105").Value = 69000
$ws.Range("L105").Value = 69000
$ws.Range("N105").Value = -75988
$ws.Range("H107").Value = 3183.2173
$ws.Range("I107").Value = 2621.875
$ws.Range("J107").Value = 4466.2856
$ws.Range("K107").Value = 2621.875
$ws.Range("L107").Value = 4466.2856
$ws.Range("M107").Value = -701.875
$ws.Range("N107").Value = -8306.285599999999
$ws.Range("H111").Value = 1210.4667
$ws.Range("I111").Value = 1046.5834
$ws.Range("K111").Value = 3139.7502
$ws.Range("M111").Value = -72.75019999999995
$ws.Range("H122").Value = 917.94116
$ws.Range("I122").Value = 760.6
$ws.Range("K122").Value = 2281.8
$ws.Range("M122").Value = 168.1999999999998
$ws.Range("H132").Value = 6853.378
$ws.Range("I132").Value = 3985.6667
$ws.Range("K132").Value = 11957.0001
$ws.Range("M132").Value = -9427.000100000001
$ws.Range("H137").Value = 5507.1665
$ws.Range("J137").Value = 2076.2727
$ws.Range("L137").Value = 6228.8181
$ws.Range("N137").Value = -11328.8181
$ws.Range("H138").Value = 22224432
$ws.Range("J138").Value = 2678.7058
$ws.Range("L138").Value = 8036.117400000001
$ws.Range("N138").Value = -18316.1174
$ws.Range("H141").Value = 8161.5386
$ws.Range("I141").Value = 10589.333
$ws.Range("J141").Value = 6080.5713
$ws.Range("K141").Value = 31767.999
$ws.Range("L141").Value = 18241.7139
$ws.Range("M141").Value = -26587.999
$ws.Range("N141").Value = -28601.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4085230.5
$ws.Range("I2").Value = 9525609
$ws.Range("J2").Value = 4946.75
$ws.Range("K2").Value = 9525609
$ws.Range("L2").Value = 4946.75
$ws.Range("M2").Value = -9525496
$ws.Range("N2").Value = -5172.75
$ws.Range("H4").Value = 1249
$ws.Range("I4").Value = 1249
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1249
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -1133
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 124
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 124
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 124
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -348
$ws.Range("H32").Value = 5285.0166
$ws.Range("I32").Value = 5615.841
$ws.Range("J32").Value = 4428.7646
$ws.Range("K32").Value = 5615.841
$ws.Range("L32").Value = 4428.7646
$ws.Range("M32").Value = -5328.841
$ws.Range("N32").Value = -5002.7646
$ws.Range("H45").Value = 8376.261
$ws.Range("I45").Value = 12444.5
$ws.Range("J45").Value = 3938.182
$ws.Range("K45").Value = 12444.5
$ws.Range("L45").Value = 3938.182
$ws.Range("M45").Value = -12067.5
$ws.Range("N45").Value = -4692.182
$ws.Range("H61").Value = 2847.9524
$ws.Range("I61").Value = 3142.4688
$ws.Range("J61").Value = 1905.5
$ws.Range("K61").Value = 3142.4688
$ws.Range("L61").Value = 1905.5
$ws.Range("M61").Value = -2930.4688
$ws.Range("N61").Value = -2329.5
$ws.Range("H74").Value = 14693.223
$ws.Range("I74").Value = 14693.223
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 14693.223
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -13819.223
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 14693.223
$ws.Range("I77").Value = 14693.223
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 73466.11500000001
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -69098.11500000001
$ws.Range("N77").ClearContents()
$ws.Range("H97").Value = 145361
$ws.Range("I97").Value = 3716.5
$ws.Range("J97").Value = 334220.34
$ws.Range("K97").Value = 3716.5
$ws.Range("L97").Value = 334220.34
$ws.Range("M97").Value = -3220.5
$ws.Range("N97").Value = -335212.34
$ws.Range("H110").Value = 3332.4285
$ws.Range("I110").Value = 2075.4
$ws.Range("K110").Value = 2075.4
$ws.Range("M110").Value = -30.40000000000009
$ws.Range("H116").Value = 4085230.5
$ws.Range("I116").Value = 9525609
$ws.Range("J116").Value = 4946.75
$ws.Range("K116").Value = 9525609
$ws.Range("L116").Value = 4946.75
$ws.Range("M116").Value = -9523315
$ws.Range("N116").Value = -9534.75
$ws.Range("H122").Value = 6050.5557
$ws.Range("I122").Value = 6053.5884
$ws.Range("K122").Value = 18160.7652
$ws.Range("M122").Value = -15710.7652
$ws.Range("H132").Value = 1623.1482
$ws.Range("I132").Value = 1531.7307
$ws.Range("K132").Value = 4595.1921
$ws.Range("M132").Value = -2065.1921
$ws.Range("H136").Value = 2847.9524
$ws.Range("I136").Value = 3142.4688
$ws.Range("J136").Value = 1905.5
$ws.Range("K136").Value = 9427.4064
$ws.Range("L136").Value = 5716.5
$ws.Range("M136").Value = -6877.4064
$ws.Range("N136").Value = -10816.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4085230.5
$ws.Range("I3").Value = 9525609
$ws.Range("J3").Value = 4946.75
$ws.Range("K3").Value = 9525609
$ws.Range("L3").Value = 4946.75
$ws.Range("M3").Value = -9525495
$ws.Range("N3").Value = -5174.75
$ws.Range("H4").Value = 124
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 124
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 124
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -354
$ws.Range("H20").Value = 1950.7222
$ws.Range("I20").Value = 1770.5814
$ws.Range("K20").Value = 1770.5814
$ws.Range("M20").Value = -1523.5814
$ws.Range("H22").Value = 850.4286
$ws.Range("I22").Value = 721.5
$ws.Range("K22").Value = 721.5
$ws.Range("M22").Value = -548.5
$ws.Range("H86").Value = 3832.2
$ws.Range("I86").Value = 3708.5
$ws.Range("J86").Value = 4079.6
$ws.Range("K86").Value = 3708.5
$ws.Range("L86").Value = 4079.6
$ws.Range("M86").Value = -2585.5
$ws.Range("N86").Value = -6325.6
$ws.Range("H89").Value = 3832.2
$ws.Range("I89").Value = 3708.5
$ws.Range("J89").Value = 4079.6
$ws.Range("K89").Value = 18542.5
$ws.Range("L89").Value = 20398
$ws.Range("M89").Value = -12926.5
$ws.Range("N89").Value = -31630
$ws.Range("H94").Value = 678.61536
$ws.Range("I94").Value = 523.7826
$ws.Range("J94").Value = 1865.6666
$ws.Range("K94").Value = 523.7826
$ws.Range("L94").Value = 1865.6666
$ws.Range("M94").Value = -72.7826
$ws.Range("N94").Value = -2767.6666
$ws.Range("H105").Value = 3926.889
$ws.Range("I105").Value = 2899.303
$ws.Range("J105").Value = 15230.333
$ws.Range("K105").Value = 2899.303
$ws.Range("L105").Value = 15230.333
$ws.Range("M105").Value = -1152.303
$ws.Range("N105").Value = -18724.333
$ws.Range("H107").Value = 1588.6
$ws.Range("I107").Value = 1248.8182
$ws.Range("J107").Value = 2163.6155
$ws.Range("K107").Value = 1248.8182
$ws.Range("L107").Value = 2163.6155
$ws.Range("M107").Value = 671.1818000000001
$ws.Range("N107").Value = -6003.6155
$ws.Range("H108").Value = 49999.9
$ws.Range("J108").Value = 49999.9
$ws.Range("L108").Value = 49999.9
$ws.Range("N108").Value = -57679.9
$ws.Range("H132").Value = 156987.19
$ws.Range("J132").Value = 156987.19
$ws.Range("L132").Value = 156987.19
$ws.Range("N132").Value = -167107.19
$ws.Range("H134").Value = 2677.3125
$ws.Range("I134").Value = 2481.9756
$ws.Range("J134").Value = 3821.4285
$ws.Range("K134").Value = 7445.926800000001
$ws.Range("L134").Value = 11464.2855
$ws.Range("M134").Value = -4910.926800000001
$ws.Range("N134").Value = -16534.2855
$ws.Range("H140").Value = 70859.89
$ws.Range("J140").Value = 70859.89
$ws.Range("L140").Value = 70859.89
$ws.Range("N140").Value = -81219.89

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1066
$ws.Range("I7").Value = 1066
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1066
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -953
$ws.Range("N7").ClearContents()
$ws.Range("H16").Value = 1934.7307
$ws.Range("I16").Value = 1558.2307
$ws.Range("K16").Value = 1558.2307
$ws.Range("M16").Value = -1271.2307
$ws.Range("H20").Value = 100386.25
$ws.Range("J20").Value = 100386.25
$ws.Range("L20").Value = 100386.25
$ws.Range("N20").Value = -100858.25
$ws.Range("H30").Value = 100386.25
$ws.Range("J30").Value = 100386.25
$ws.Range("L30").Value = 100386.25
$ws.Range("N30").Value = -100568.25
$ws.Range("H31").Value = 2628.2354
$ws.Range("J31").Value = 2937.6667
$ws.Range("L31").Value = 2937.6667
$ws.Range("N31").Value = -3527.6667
$ws.Range("H34").Value = 2628.2354
$ws.Range("J34").Value = 2937.6667
$ws.Range("L34").Value = 2937.6667
$ws.Range("N34").Value = -3341.6667
$ws.Range("H58").Value = 1320.6666
$ws.Range("I58").Value = 1306.9333
$ws.Range("K58").Value = 1306.9333
$ws.Range("M58").Value = -1103.9333
$ws.Range("H64").Value = 59999.668
$ws.Range("J64").Value = 59999.668
$ws.Range("L64").Value = 59999.668
$ws.Range("N64").Value = -60495.668
$ws.Range("H67").Value = 59999.668
$ws.Range("J67").Value = 59999.668
$ws.Range("L67").Value = 59999.668
$ws.Range("N67").Value = -61715.668
$ws.Range("H105").Value = 2592.739
$ws.Range("J105").Value = 16666.5
$ws.Range("L105").Value = 16666.5
$ws.Range("N105").Value = -20160.5
$ws.Range("H109").Value = 49771
$ws.Range("J109").Value = 49771
$ws.Range("L109").Value = 49771
$ws.Range("N109").Value = -51851
$ws.Range("H113").Value = 1934.7307
$ws.Range("I113").Value = 1558.2307
$ws.Range("K113").Value = 1558.2307
$ws.Range("M113").Value = 611.7692999999999
$ws.Range("H122").Value = 7049.577
$ws.Range("I122").Value = 9744.25
$ws.Range("J122").Value = 2738.1
$ws.Range("K122").Value = 29232.75
$ws.Range("L122").Value = 8214.299999999999
$ws.Range("M122").Value = -26782.75
$ws.Range("N122").Value = -13114.3
$ws.Range("H128").Value = 100386.25
$ws.Range("J128").Value = 100386.25
$ws.Range("L128").Value = 100386.25
$ws.Range("N128").Value = -110346.25
$ws.Range("H132").Value = 5963.814
$ws.Range("I132").Value = 2387.111
$ws.Range("J132").Value = 24358.285
$ws.Range("K132").Value = 7161.333
$ws.Range("L132").Value = 73074.855
$ws.Range("M132").Value = -4631.333
$ws.Range("N132").Value = -78134.855
$ws.Range("H134").Value = 1601.2325
$ws.Range("I134").Value = 1474
$ws.Range("K134").Value = 4422
$ws.Range("M134").Value = -1887
$ws.Range("H136").Value = 1320.6666
$ws.Range("I136").Value = 1306.9333
$ws.Range("K136").Value = 3920.7999
$ws.Range("M136").Value = -1370.7999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1499.75
$ws.Range("I22").Value = 499
$ws.Range("J22").Value = 1833.3334
$ws.Range("K22").Value = 1497
$ws.Range("L22").Value = 5500.0002
$ws.Range("M22").Value = -1328
$ws.Range("N22").Value = -5838.0002
$ws.Range("H27").Value = 1499.75
$ws.Range("I27").Value = 499
$ws.Range("J27").Value = 1833.3334
$ws.Range("K27").Value = 1497
$ws.Range("L27").Value = 5500.0002
$ws.Range("M27").Value = -1395
$ws.Range("N27").Value = -5704.0002
$ws.Range("H46").Value = 530.25
$ws.Range("I46").Value = 183
$ws.Range("J46").Value = 738.6
$ws.Range("K46").Value = 549
$ws.Range("L46").Value = 2215.8
$ws.Range("M46").Value = -458
$ws.Range("N46").Value = -2397.8
$ws.Range("H63").Value = 999.5
$ws.Range("I63").Value = 999.5
$ws.Range("K63").Value = 2998.5
$ws.Range("M63").Value = -2249.5
$ws.Range("H66").Value = 999.5
$ws.Range("I66").Value = 999.5
$ws.Range("K66").Value = 8995.5
$ws.Range("M66").Value = -5251.5
$ws.Range("H68").Value = 1224.5
$ws.Range("J68").Value = 974
$ws.Range("L68").Value = 2922
$ws.Range("N68").Value = -4544
$ws.Range("H71").Value = 1224.5
$ws.Range("J71").Value = 974
$ws.Range("L71").Value = 8766
$ws.Range("N71").Value = -16878
$ws.Range("H74").Value = 4666
$ws.Range("I74").Value = 4666
$ws.Range("K74").Value = 13998
$ws.Range("M74").Value = -12937
$ws.Range("H77").Value = 4666
$ws.Range("I77").Value = 4666
$ws.Range("K77").Value = 41994
$ws.Range("M77").Value = -36690
$ws.Range("H94").Value = 174242910
$ws.Range("I94").Value = 738.5
$ws.Range("K94").Value = 2215.5
$ws.Range("M94").Value = -1539.5
$ws.Range("H112").Value = 1406
$ws.Range("I112").Value = 1406
$ws.Range("K112").Value = 4218
$ws.Range("M112").Value = -3110
$ws.Range("H113").Value = 587.8889
$ws.Range("I113").Value = 446.2857
$ws.Range("J113").Value = 740.38464
$ws.Range("K113").Value = 1338.8571
$ws.Range("L113").Value = 2221.15392
$ws.Range("M113").Value = 831.1428999999998
$ws.Range("N113").Value = -6561.15392
$ws.Range("H117").Value = 3582.5
$ws.Range("J117").Value = 4750
$ws.Range("L117").Value = 14250
$ws.Range("N117").Value = -21134
$ws.Range("H121").Value = 1575.5
$ws.Range("J121").Value = 1587.7273
$ws.Range("L121").Value = 4763.1819
$ws.Range("N121").Value = -7383.1819
$ws.Range("H129").Value = 1190.2941
$ws.Range("J129").Value = 1663
$ws.Range("L129").Value = 4989
$ws.Range("N129").Value = -14989
$ws.Range("H140").Value = 1921.7778
$ws.Range("I140").Value = 1921.7778
$ws.Range("K140").Value = 5765.3334
$ws.Range("M140").Value = -585.3334000000004
$ws.Range("H141").Value = 3034.8
$ws.Range("I141").Value = 3034.8
$ws.Range("K141").Value = 9104.400000000001
$ws.Range("M141").Value = -3924.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9999.333000000001
$ws.Range("J70").Value = 10000
$ws.Range("L70").Value = 10000
$ws.Range("N70").Value = -10540
$ws.Range("H73").Value = 9999.333000000001
$ws.Range("J73").Value = 10000
$ws.Range("L73").Value = 10000
$ws.Range("N73").Value = -11872
$ws.Range("H80").Value = 76669460
$ws.Range("I80").Value = 115002310
$ws.Range("J80").Value = 3749.5
$ws.Range("K80").Value = 115002310
$ws.Range("L80").Value = 3749.5
$ws.Range("M80").Value = -115001312
$ws.Range("N80").Value = -5745.5
$ws.Range("H83").Value = 76669460
$ws.Range("I83").Value = 115002310
$ws.Range("J83").Value = 3749.5
$ws.Range("K83").Value = 575011550
$ws.Range("L83").Value = 18747.5
$ws.Range("M83").Value = -575006558
$ws.Range("N83").Value = -28731.5
$ws.Range("H92").Value = 29999
$ws.Range("J92").Value = 29999
$ws.Range("L92").Value = 29999
$ws.Range("N92").Value = -33743
$ws.Range("H102").Value = 4218.5
$ws.Range("I102").Value = 3736.9412
$ws.Range("J102").Value = 5855.8
$ws.Range("K102").Value = 3736.9412
$ws.Range("L102").Value = 5855.8
$ws.Range("M102").Value = -2114.9412
$ws.Range("N102").Value = -9099.799999999999
$ws.Range("H132").Value = 4202.654
$ws.Range("I132").Value = 3969.5417
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 11908.6251
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -9378.625100000001
$ws.Range("N132").Value = -26060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6777.5
$ws.Range("I7").Value = 5535
$ws.Range("K7").Value = 5535
$ws.Range("M7").Value = -5423
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H22").Value = 771.3333
$ws.Range("I22").Value = 699
$ws.Range("K22").Value = 699
$ws.Range("M22").Value = -404
$ws.Range("H27").Value = 771.3333
$ws.Range("I27").Value = 699
$ws.Range("K27").Value = 699
$ws.Range("M27").Value = -592
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H46").Value = 2885.1365
$ws.Range("I46").Value = 1929
$ws.Range("J46").Value = 3841.2727
$ws.Range("K46").Value = 1929
$ws.Range("L46").Value = 3841.2727
$ws.Range("M46").Value = -1741
$ws.Range("N46").Value = -4217.2727
$ws.Range("H55").Value = 413.42105
$ws.Range("I55").Value = 173.73334
$ws.Range("K55").Value = 173.73334
$ws.Range("M55").Value = -0.7333399999999983
$ws.Range("H61").Value = 3251.5
$ws.Range("I61").Value = 5497.5
$ws.Range("J61").Value = 1005.5
$ws.Range("K61").Value = 5497.5
$ws.Range("L61").Value = 1005.5
$ws.Range("M61").Value = -5295.5
$ws.Range("N61").Value = -1409.5
$ws.Range("H68").Value = 4168.1665
$ws.Range("I68").Value = 849.92
$ws.Range("J68").Value = 20759.4
$ws.Range("K68").Value = 849.92
$ws.Range("L68").Value = 20759.4
$ws.Range("M68").Value = -100.92
$ws.Range("N68").Value = -22257.4
$ws.Range("H71").Value = 4168.1665
$ws.Range("I71").Value = 849.92
$ws.Range("J71").Value = 20759.4
$ws.Range("K71").Value = 4249.599999999999
$ws.Range("L71").Value = 103797
$ws.Range("M71").Value = -505.5999999999995
$ws.Range("N71").Value = -111285
$ws.Range("H82").Value = 1977.99
$ws.Range("I82").Value = 1942.7084
$ws.Range("J82").Value = 2824.75
$ws.Range("K82").Value = 1942.7084
$ws.Range("L82").Value = 2824.75
$ws.Range("M82").Value = -1581.7084
$ws.Range("N82").Value = -3546.75
$ws.Range("H85").Value = 1977.99
$ws.Range("I85").Value = 1942.7084
$ws.Range("J85").Value = 2824.75
$ws.Range("K85").Value = 1942.7084
$ws.Range("L85").Value = 2824.75
$ws.Range("M85").Value = -694.7084
$ws.Range("N85").Value = -5320.75
$ws.Range("H93").Value = 2001.7894
$ws.Range("I93").Value = 1835.2222
$ws.Range("J93").Value = 5000
$ws.Range("K93").Value = 1835.2222
$ws.Range("L93").Value = 5000
$ws.Range("M93").Value = -587.2221999999999
$ws.Range("N93").Value = -7496
$ws.Range("H94").Value = 40000
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41352
$ws.Range("H100").Value = 3332.5
$ws.Range("I100").Value = 3665
$ws.Range("K100").Value = 3665
$ws.Range("M100").Value = -3124
$ws.Range("H113").Value = 3251.5
$ws.Range("I113").Value = 5497.5
$ws.Range("J113").Value = 1005.5
$ws.Range("K113").Value = 5497.5
$ws.Range("L113").Value = 1005.5
$ws.Range("M113").Value = -3327.5
$ws.Range("N113").Value = -5345.5
$ws.Range("H123").Value = 50933
$ws.Range("J123").Value = 50933
$ws.Range("L123").Value = 50933
$ws.Range("N123").Value = -60733
$ws.Range("H126").Value = 6777.5
$ws.Range("I126").Value = 5535
$ws.Range("K126").Value = 16605
$ws.Range("M126").Value = -14135
$ws.Range("H132").Value = 36804.473
$ws.Range("I132").Value = 40610.59
$ws.Range("J132").Value = 4452.5
$ws.Range("K132").Value = 121831.77
$ws.Range("L132").Value = 13357.5
$ws.Range("M132").Value = -119301.77
$ws.Range("N132").Value = -18417.5
$ws.Range("H136").Value = 2362.4
$ws.Range("I136").Value = 1571.5
$ws.Range("J136").Value = 10798.667
$ws.Range("K136").Value = 4714.5
$ws.Range("L136").Value = 32396.001
$ws.Range("M136").Value = -2164.5
$ws.Range("N136").Value = -37496.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 8333
$ws.Range("J30").Value = 8333
$ws.Range("L30").Value = 8333
$ws.Range("N30").Value = -8547
$ws.Range("H62").Value = 15470.429
$ws.Range("I62").Value = 11459
$ws.Range("J62").Value = 25499
$ws.Range("K62").Value = 11459
$ws.Range("L62").Value = 25499
$ws.Range("M62").Value = -10835
$ws.Range("N62").Value = -26747
$ws.Range("H65").Value = 15470.429
$ws.Range("I65").Value = 11459
$ws.Range("J65").Value = 25499
$ws.Range("K65").Value = 57295
$ws.Range("L65").Value = 127495
$ws.Range("M65").Value = -54175
$ws.Range("N65").Value = -133735
$ws.Range("H92").Value = 36999.75
$ws.Range("J92").Value = 36999.75
$ws.Range("L92").Value = 36999.75
$ws.Range("N92").Value = -41991.75
$ws.Range("H107").Value = 4366.3335
$ws.Range("I107").Value = 4366.3335
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 13099.0005
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -11179.0005
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 2178.9092
$ws.Range("J113").Value = 3497.6667
$ws.Range("L113").Value = 10493.0001
$ws.Range("N113").Value = -14833.0001
$ws.Range("H132").Value = 3358.45
$ws.Range("I132").Value = 2897.7144
$ws.Range("J132").Value = 6583.6
$ws.Range("K132").Value = 8693.143199999999
$ws.Range("L132").Value = 19750.8
$ws.Range("M132").Value = -6163.143199999999
$ws.Range("N132").Value = -24810.8
$ws.Range("H135").Value = 241940.2
$ws.Range("J135").Value = 241940.2
$ws.Range("L135").Value = 241940.2
$ws.Range("N135").Value = -252080.2
$ws.Range("H136").Value = 2556.4783
$ws.Range("I136").Value = 2671.3809
$ws.Range("J136").Value = 1350
$ws.Range("K136").Value = 8014.1427
$ws.Range("L136").Value = 4050
$ws.Range("M136").Value = -5464.1427
$ws.Range("N136").Value = -9150

